$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 89.2
$ws.Range("I2").Value = 103
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 103
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = -306
$ws.Range("H8").Value = 700
$ws.Range("I8").Value = 700
$ws.Range("K8").Value = 2100
$ws.Range("M8").Value = -1961
$ws.Range("H51").Value = 6874
$ws.Range("I51").Value = 9330.333000000001
$ws.Range("J51").Value = 5821.2856
$ws.Range("K51").Value = 9330.333000000001
$ws.Range("L51").Value = 5821.2856
$ws.Range("M51").Value = -8846.333000000001
$ws.Range("N51").Value = -6789.2856
$ws.Range("H53").Value = 72.666664
$ws.Range("I53").Value = 70.27273
$ws.Range("J53").Value = 99
$ws.Range("K53").Value = 70.27273
$ws.Range("L53").Value = 99
$ws.Range("M53").Value = 566.72727
$ws.Range("N53").Value = -1373
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142
$ws.Range("H80").Value = 485.5
$ws.Range("I80").Value = 383.33334
$ws.Range("K80").Value = 1150.00002
$ws.Range("M80").Value = -152.0000199999999
$ws.Range("H83").Value = 485.5
$ws.Range("I83").Value = 383.33334
$ws.Range("K83").Value = 3450.00006
$ws.Range("M83").Value = 1541.99994
$ws.Range("H92").Value = 1249.6666
$ws.Range("I92").Value = 1100
$ws.Range("K92").Value = 1100
$ws.Range("M92").Value = 148
$ws.Range("H118").Value = 835.7143
$ws.Range("I118").Value = 337.5
$ws.Range("K118").Value = 1012.5
$ws.Range("M118").Value = 644.5
$ws.Range("H129").Value = 1446.1666
$ws.Range("I129").Value = 1446.1666
$ws.Range("K129").Value = 4338.4998
$ws.Range("M129").Value = 661.5002000000004
$ws.Range("H132").Value = 3076.8
$ws.Range("I132").Value = 1135.3334
$ws.Range("K132").Value = 3406.0002
$ws.Range("M132").Value = -876.0001999999999
$ws.Range("H133").Value = 99779.664
$ws.Range("J133").Value = 99779.664
$ws.Range("L133").Value = 99779.664
$ws.Range("N133").Value = -109899.664
$ws.Range("H137").Value = 2664.7307
$ws.Range("I137").Value = 1920.3334
$ws.Range("J137").Value = 3058.8235
$ws.Range("K137").Value = 5761.0002
$ws.Range("L137").Value = 9176.470499999999
$ws.Range("M137").Value = -3211.0002
$ws.Range("N137").Value = -14276.4705
$ws.Range("H138").Value = 3112.138
$ws.Range("I138").Value = 784.125
$ws.Range("K138").Value = 2352.375
$ws.Range("M138").Value = 2787.625
$ws.Range("H141").Value = 49998.75
$ws.Range("I141").Value = 49998.75
$ws.Range("K141").Value = 149996.25
$ws.Range("M141").Value = -144816.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3704.238
$ws.Range("I32").Value = 3864.45
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 3864.45
$ws.Range("L32").Value = 500
$ws.Range("M32").Value = -3577.45
$ws.Range("N32").Value = -1074
$ws.Range("H45").Value = 2640.9473
$ws.Range("I45").Value = 2475.4707
$ws.Range("K45").Value = 2475.4707
$ws.Range("M45").Value = -2098.4707
$ws.Range("H97").Value = 1022.4
$ws.Range("I97").Value = 697
$ws.Range("K97").Value = 697
$ws.Range("M97").Value = -201
$ws.Range("H110").Value = 2996.6155
$ws.Range("I110").Value = 492.83334
$ws.Range("K110").Value = 492.83334
$ws.Range("M110").Value = 1552.16666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2756.625
$ws.Range("I94").Value = 2508.8333
$ws.Range("J94").Value = 3500
$ws.Range("K94").Value = 2508.8333
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = -2057.8333
$ws.Range("N94").Value = -4402
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12499
$ws.Range("I31").Value = 12499
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 12499
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -12204
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 12499
$ws.Range("I34").Value = 12499
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 12499
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -12297
$ws.Range("N34").ClearContents()
$ws.Range("H86").Value = 5399.8
$ws.Range("I86").Value = 4749.75
$ws.Range("K86").Value = 4749.75
$ws.Range("M86").Value = -3626.75
$ws.Range("H89").Value = 5399.8
$ws.Range("I89").Value = 4749.75
$ws.Range("K89").Value = 23748.75
$ws.Range("M89").Value = -18132.75
$ws.Range("H107").Value = 179.4
$ws.Range("I107").Value = 99.333336
$ws.Range("J107").Value = 299.5
$ws.Range("K107").Value = 99.333336
$ws.Range("L107").Value = 299.5
$ws.Range("M107").Value = 1820.666664
$ws.Range("N107").Value = -4139.5
$ws.Range("H134").Value = 1404.375
$ws.Range("I134").Value = 1539.1666
$ws.Range("K134").Value = 4617.4998
$ws.Range("M134").Value = -2082.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 999999
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H26").Value = 75
$ws.Range("I26").Value = 75
$ws.Range("K26").Value = 225
$ws.Range("M26").Value = 63
$ws.Range("H40").Value = 265.69232
$ws.Range("I40").Value = 182.66667
$ws.Range("K40").Value = 730.66668
$ws.Range("M40").Value = -661.66668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7214.5884
$ws.Range("I122").Value = 7415.5
$ws.Range("K122").Value = 22246.5
$ws.Range("M122").Value = -19796.5
$ws.Range("H126").Value = 1688.4
$ws.Range("I126").Value = 1661.5
$ws.Range("K126").Value = 4984.5
$ws.Range("M126").Value = -2514.5
$ws.Range("H132").Value = 3145.875
$ws.Range("I132").Value = 3145.875
$ws.Range("K132").Value = 9437.625
$ws.Range("M132").Value = -6907.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 860
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 933.3333
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 933.3333
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -1523.3333
$ws.Range("H27").Value = 860
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 933.3333
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 933.3333
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -1147.3333
$ws.Range("H46").Value = 959.8
$ws.Range("I46").Value = 866.6667
$ws.Range("K46").Value = 866.6667
$ws.Range("M46").Value = -678.6667
$ws.Range("H55").Value = 1026.0769
$ws.Range("I55").Value = 268.125
$ws.Range("K55").Value = 268.125
$ws.Range("M55").Value = -95.125
$ws.Range("H93").Value = 1762.5
$ws.Range("I93").Value = 2016.6666
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 2016.6666
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -768.6666
$ws.Range("N93").Value = -3496
$ws.Range("H100").Value = 4399
$ws.Range("I100").Value = 4997.5
$ws.Range("K100").Value = 4997.5
$ws.Range("M100").Value = -4456.5
$ws.Range("H122").Value = 4716.143
$ws.Range("I122").Value = 5201.6
$ws.Range("J122").Value = 3502.5
$ws.Range("K122").Value = 15604.8
$ws.Range("L122").Value = 10507.5
$ws.Range("M122").Value = -13154.8
$ws.Range("N122").Value = -15407.5
$ws.Range("H132").Value = 9829
$ws.Range("I132").Value = 9829
$ws.Range("K132").Value = 29487
$ws.Range("M132").Value = -26957
$ws.Range("H136").Value = 5427.375
$ws.Range("I136").Value = 5532.143
$ws.Range("J136").Value = 4694
$ws.Range("K136").Value = 16596.429
$ws.Range("L136").Value = 14082
$ws.Range("M136").Value = -14046.429
$ws.Range("N136").Value = -19182

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 456.14285
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 4210
$ws.Range("I132").Value = 4210
$ws.Range("K132").Value = 12630
$ws.Range("M132").Value = -10100
